$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "60.727.23"
    "E2"  = "  -0.25%  "
    "D3"  = "2.396.64"
    "E3"  = "  -0.86%  "
    "E4"  = "  +0.45%  "
    "D5"  = "560.64"
    "E5"  = "  -1.54%  "
    "D6"  = "141.30"
    "E6"  = "  +1.22%  "
    "E7"  = "  -0.33%  "
    "D8"  = "0.539"
    "E8"  = "  +2.59%  "
    "D9"  = "2.402.48"
    "E9"  = "  +0.05%  "
    "E10" = "  +0.69%  "
    "E11" = "  -0.29%  "
    "E12" = "  +1.85%  "
    "E13" = "  +2.18%  "
    "D14" = "26.30"
    "E14" = "  +0.81%  "
    "E15" = "  -1.23%  "
    "D16" = "2.777.69"
    "E16" = "  -1.82%  "
    "D17" = "60.319.69"
    "E17" = "  -0.88%  "
    "D18" = "2.401.18"
    "E18" = "  -0.18%  "
    "D19" = "8.13"
    "E19" = "  +6.13%  "
    "D20" = "10.67"
    "E20" = "  +0.36%  "
    "D21" = "323.92"
    "E21" = "  +0.42%  "
    "E22" = "  +1.36%  "
    "E23" = "  -0.40%  "
    "E24" = "  -0.20%  "
    "E25" = "  +0.96%  "
    "D26" = "64.73"
    "E26" = "  -0.06%  "
    "D27" = "572.80"
    "E27" = "  -1.20%  "
    "E28" = "  -2.54%  "
    "D29" = "2.513.79"
    "E29" = "  -1.06%  "
    "D30" = "0.0₃0935"
    "E30" = "  +0.61%  "
    "D31" = "8.04"
    "E31" = "  +2.31%  "
    "D32" = "1.34"
    "E32" = "  -0.35%  "
    "E33" = "  -1.71%  "
    "E34" = "  -0.16%  "
    "E35" = "  -0.56%  "
    "E36" = "  +3.87%  "
    "D37" = "152.04"
    "E37" = "  -0.11%  "
    "E38" = "  +0.93%  "
    "E39" = "  +0.06%  "
    "D40" = "18.27"
    "E40" = "  +0.20%  "
    "D41" = "5.17"
    "E41" = "  +0.87%  "
    "E42" = "  -0.08%  "
    "D43" = "2.51"
    "E43" = "  +6.89%  "
    "E44" = "  +0.74%  "
    "E45" = "  +1.14%  "
    "D46" = "0.0₆0281"
    "E46" = "  +2.12%  "
    "D47" = "141.54"
    "E47" = "  -0.32%  "
    "E48" = "  +0.32%  "
    "E49" = "  +0.28%  "
    "D50" = "0.0507"
    "E50" = "  +1.25%  "
    "D51" = "19.39"
    "E51" = "  +0.19%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    # Force text format so numeric-looking values (e.g. "141.30", "560.64")
    # keep their exact string representation instead of becoming numbers,
    # then restore the original cell style so formatting is unaffected.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = $origStyle
}
